$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("task")
$ws.Activate()

# Remove the "Reviewed by" (column L) entries for rows 7-10 (they are now blank)
$ws.Range("L7:L10").ClearContents()

# Add a new row 11 with a new task
$ws.Range("B11").Value = "Add view history for a vahicle in service form"
$ws.Range("D11").Value = "Functional"
$ws.Range("F11").Value = "Ruwan"
$ws.Range("H11").Value = "1 day"
$ws.Range("J11").Value = "new"
$ws.Range("N10").Copy($ws.Range("N11"))
$ws.Range("N11").Value = 42987

$ws.Range("P10").Copy($ws.Range("P11"))
$ws.Range("P11").Value = 42982

# Update the selection to match the new active cell
$ws.Range("I18").Select()
